$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "Context (Metaclass): SubjectKind. Employment" "Context (Metaclass): SubjectKind. Works, Employs, Performs"
Replace-Text "Subject (Context): AggregatedReifiedAggregation. Employer, Employee, Position" "Subject (Context): AggregatedReifiedAggregation SKs. HasEmployer, HasEmployee, HasPosition"
Replace-Text "Predicate (Role): AggregationSubjectKind. Employer / Employee / Position" "Predicate (Role): AggregationSubjectKinds. Employer / Employee / Position"
Replace-Text "Object (Occurrence): AggregatedReifiedAggregation. Employer, Employee, Position" "Object (Occurrence): AggregatedReifiedAggregation OKs. EmployerOf, EmployeeOf, PositionOf"
Replace-Text "Context (Dimension): SubjectKind. LaboralStatus" "Context (Dimension): SubjectKind. Employment"
Replace-Text "Subject (Measure): AggregatedReifiedActivation. HasEmployer" "Subject (Measure): AggregatedReifiedActivation SKs. Employed, Employing, Performing"
Replace-Text "Predicate (Unit): ActivationSubjectKind. Employment" "Predicate (Unit): ActivationSubjectKinds. Works, Employs, Performs"
Replace-Text "Object (Value): AggregatedReifiedActivation. HasPosition" "Object (Value): AggregatedReifiedActivation OKs. EmployeedAt, Employing, PerformingPosition"
Replace-Text "Activation (S, O) from Activation Predicate / Aggregation Subject Kind Attributes (PKs SK / OK)." "Activation (S, O) from Activation Predicate / Aggregation Subject Kind Attributes (PKs SKs / OKs)."
Replace-Text "Dimensional (S, O) from Alignment Predicate / Activation Subject Kind Attributes (PKs SK / OK)." "Dimensional (S, O) from Alignment Predicate / Activation Subject Kind Attributes (PKs SKs / OKs)."
